$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$data = @(
    @("2023-07-17 02:02:19", "Bank", "stress", 10, 8, 2, 0.004945099992587836),
    @("2023-07-17 02:03:18", "Bank", "stress", 10, 8, 2, 0.004945099992587836),
    @("2023-07-17 02:03:18", "Bank", "stress", 7, 5, 2, 0.01952064496563264),
    @("2023-07-17 02:03:18", "Bank", "stress", 2, 1, 1, 0.5680354784012266),
    @("2023-07-17 02:03:18", "Bank", "stress", 6, 5, 1, 0.02760636085666775),
    @("2023-07-17 02:03:19", "Bank", "stress", 5, 4, 1, 0.04585823797499209),
    @("2023-07-17 02:03:19", "Bank", "stress", 3, 2, 1, 0.2352084552295149),
    @("2023-07-17 02:03:19", "Bank", "stress", 8, 7, 1, 0.01383268926919492),
    @("2023-07-17 02:03:19", "Bank", "stress", 1, 0, 1, 1)
)

$startRow = 16
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]
    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
    $ws.Cells.Item($row, 4).Value = $vals[3]
    $ws.Cells.Item($row, 5).Value = $vals[4]
    $ws.Cells.Item($row, 6).Value = $vals[5]
    $ws.Cells.Item($row, 7).Value = $vals[6]
}
